$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 369, shifting all existing rows (369-459) down to (370-460)
$ws.Rows(369).Insert()

# Populate the newly inserted row 369 with the new weekly price record
$ws.Range("A369").Value = 8
$ws.Range("B369").Value = "Terminal La Palmera de La Serena"
$ws.Range("C369").Value = "Coquimbo"
$ws.Range("D369").Value = 44782
$ws.Range("E369").Value = 4
$ws.Range("F369").Value = 100114001
$ws.Range("G369").Value = "Papa"
$ws.Range("H369").Value = "Asterix"
$ws.Range("I369").Value = "1a (cosecha)"
$ws.Range("J369").Value = 2000
$ws.Range("K369").Value = 11000
$ws.Range("L369").Value = 12000
$ws.Range("M369").Value = 11500
$ws.Range("N369").Value = "$/saco 25 kilos"
$ws.Range("O369").Value = "Provincia de Melipilla"
$ws.Range("P369").Value = 460
$ws.Range("Q369").Value = 25
$ws.Range("R369").Value = "Hortaliza"
